# Auto-generated edit script to update F-column values
# ("想去人数" - number of people interested) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6576
$ws.Range("F3").Value = 766
$ws.Range("F4").Value = 1101
$ws.Range("F6").Value = 625
$ws.Range("F7").Value = 213
$ws.Range("F8").Value = 44
$ws.Range("F9").Value = 800
$ws.Range("F10").Value = 1264
$ws.Range("F11").Value = 24
$ws.Range("F13").Value = 515
$ws.Range("F14").Value = 504
$ws.Range("F16").Value = 1042
$ws.Range("F17").Value = 1453
$ws.Range("F19").Value = 427
$ws.Range("F20").Value = 428
$ws.Range("F22").Value = 1092
$ws.Range("F23").Value = 206
$ws.Range("F24").Value = 2284
$ws.Range("F26").Value = 152
$ws.Range("F27").Value = 419
$ws.Range("F29").Value = 3720

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 32
$ws.Range("F6").Value = 731
$ws.Range("F9").Value = 1028
$ws.Range("F11").Value = 132
$ws.Range("F14").Value = 5
$ws.Range("F17").Value = 387
$ws.Range("F25").Value = 235
$ws.Range("F32").Value = 1690

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1224
$ws.Range("F5").Value = 1602
$ws.Range("F8").Value = 908

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1224
$ws.Range("F4").Value = 1602
$ws.Range("F7").Value = 908
$ws.Range("F8").Value = 6576
$ws.Range("F9").Value = 32
$ws.Range("F10").Value = 766
$ws.Range("F11").Value = 731
$ws.Range("F13").Value = 625
$ws.Range("F14").Value = 213
$ws.Range("F15").Value = 44
$ws.Range("F16").Value = 800
$ws.Range("F19").Value = 132
$ws.Range("F20").Value = 132
$ws.Range("F23").Value = 1264
$ws.Range("F24").Value = 24
$ws.Range("F26").Value = 515
$ws.Range("F27").Value = 504
$ws.Range("F29").Value = 387
$ws.Range("F32").Value = 1042
$ws.Range("F33").Value = 1453
$ws.Range("F36").Value = 427
$ws.Range("F37").Value = 428
$ws.Range("F40").Value = 1092
$ws.Range("F41").Value = 206
$ws.Range("F42").Value = 2284
$ws.Range("F43").Value = 1690
$ws.Range("F44").Value = 1690
$ws.Range("F45").Value = 152
$ws.Range("F46").Value = 419
$ws.Range("F47").Value = 3720
